$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Insert a new "CreditAmt" column right after "DebitAmt" (col F), which
# shifts the LowWait..HighestExplicitWait columns from G:N to H:O.
$ws.Columns("G").Insert()
$ws.Range("G1").Value = "CreditAmt"

# New account code + matching Debit/Credit amounts for the sample row.
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Value = "260.38533.230804.40328.100.00000.000.0000.0000"
$ws.Range("F2").Value = "275"
$ws.Range("G2").Value = "275"

# Update the instructional note and extend/merge it across the new column.
$ws.Range("F5").Value = "Credit Amt and Debit Amt should be same"
$ws.Range("F5:G5").Merge()
